$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray formatted-but-empty row 15 (G15 had a fill-only style, no value)
$ws.Rows(15).Delete()

# Add the new respondent row (row 11): "Basanta Shrestha"
$ws.Range("A11").Value = "Basanta Shrestha"

# Column B: Open Science (same formatting as E2)
$ws.Range("E2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Open Science"

# Column C: Reproducibility crisis (same formatting as C2)
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "Reproducibility crisis"

# Column D: AI (same formatting as E3)
$ws.Range("E3").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "AI"

# Column E: Coding Practices (same formatting as D3)
$ws.Range("D3").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "Coding Practices"

# Column F: Citizen Science (same formatting as H2)
$ws.Range("H2").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = "Citizen Science"

# Column G: Modern Computing (same formatting as G4)
$ws.Range("G4").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G11").Value = "Modern Computing"

# Column H: p-value (new term, reuses the "p value conundrum" formatting from G7)
$ws.Range("G7").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = "p-value"

# Column I: Data Repositories (same formatting as K3)
$ws.Range("K3").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = "Data Repositories"

# Column J: Predatory Journals (same formatting as I2)
$ws.Range("I2").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("J11").Value = "Predatory Journals"

# Column K: Pay to publish (same formatting as K2)
$ws.Range("K2").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = "Pay to publish"

# Column L: Mental Health (same formatting as L2)
$ws.Range("L2").Copy()
$ws.Range("L11").PasteSpecial(-4122)
$ws.Range("L11").Value = "Mental Health"

# Column M: Job Application Packets (same formatting as M2)
$ws.Range("M2").Copy()
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("M11").Value = "Job Application Packets"

# Column N: Presentation Woes (same formatting as N3)
$ws.Range("N3").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = "Presentation Woes"

# Column O: Keeping up with Scientific Literature (same formatting as O3)
$ws.Range("O3").Copy()
$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = "Keeping up with Scientific Literature"

$excel.CutCopyMode = $false

# Restore the selection to where the author left off
$ws.Range("I15").Select()
